$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in "Day 6" (column G) sleep-diary entries for the last week's table
# (rows 122-135), which were previously left blank.
$ws.Range("G122").Value = "7：22"
$ws.Range("G123").Value = "7：30"
$ws.Range("G124").Value = "23：20"
$ws.Range("G125").Value = "23：30"
$ws.Range("G126").Value = 5
$ws.Range("G127").Value = 0
$ws.Range("G128").Value = 0
$ws.Range("G129").Value = 480
$ws.Range("G130").Value = "无"
$ws.Range("G131").Value = "无"
$ws.Range("G132").Value = 4
$ws.Range("G133").Value = 4
$ws.Range("G134").Value = 4
$ws.Range("G135").Value = "无"

# Move the active selection to match where the author ended up.
$ws.Range("G135").Select()
